$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest")

$ws.Range("B2").Value = 0.5085426078324898
$ws.Range("C2").Value = 0.6312770643337786
$ws.Range("D2").Value = 0.5775225091701024
$ws.Range("E2").Value = 0.7599490174808455
$ws.Range("F2").Value = 0.5845391272395405
$ws.Range("G2").Value = 15
$ws.Range("B3").Value = 0.3604583747834282
$ws.Range("C3").Value = 0.5284892497187151
$ws.Range("D3").Value = 0.3749043735392947
$ws.Range("E3").Value = 0.6122943520393559
$ws.Range("F3").Value = 0.5136324539103049
$ws.Range("G3").Value = 14
$ws.Range("B4").Value = 0.2671644164584084
$ws.Range("C4").Value = 0.4711099050880848
$ws.Range("D4").Value = 0.3212841183425387
$ws.Range("E4").Value = 0.5668192995501641
$ws.Range("F4").Value = 0.5203199983321719
$ws.Range("G4").Value = 13
$ws.Range("B5").Value = 0.4282641160228051
$ws.Range("C5").Value = 0.548481490052795
$ws.Range("D5").Value = 0.4014706193385897
$ws.Range("E5").Value = 0.6336170920505457
$ws.Range("F5").Value = 0.4877336824715215
$ws.Range("G5").Value = 12
$ws.Range("B6").Value = 0.4297788858055521
$ws.Range("C6").Value = 0.5712726212628527
$ws.Range("D6").Value = 0.4264280183742889
$ws.Range("E6").Value = 0.6530145621456607
$ws.Range("F6").Value = 0.5156451691415619
$ws.Range("G6").Value = 11
$ws.Range("B7").Value = 0.3783814472866451
$ws.Range("C7").Value = 0.5371069045676157
$ws.Range("D7").Value = 0.3902711452809265
$ws.Range("E7").Value = 0.6247168520865485
$ws.Range("F7").Value = 0.5239790343878171
$ws.Range("G7").Value = 10
$ws.Range("B8").Value = 0.3326355726653664
$ws.Range("C8").Value = 0.4945997455591428
$ws.Range("D8").Value = 0.3233547881867544
$ws.Range("E8").Value = 0.5686429355815075
$ws.Range("F8").Value = 0.4891798334788348
$ws.Range("G8").Value = 9
$ws.Range("B9").Value = 0.3620304170102688
$ws.Range("C9").Value = 0.5407351962130211
$ws.Range("D9").Value = 0.3743586110673328
$ws.Range("E9").Value = 0.6118485197067431
$ws.Range("F9").Value = 0.5273032071390139
$ws.Range("G9").Value = 8
$ws.Range("B10").Value = 0.4970330926187166
$ws.Range("C10").Value = 0.5512696337016082
$ws.Range("D10").Value = 0.4018017306810739
$ws.Range("E10").Value = 0.6338783248235216
$ws.Range("F10").Value = 0.4249154521118755
$ws.Range("G10").Value = 7
$ws.Range("B11").Value = 0.4078060000532727
$ws.Range("C11").Value = 0.4615113944401893
$ws.Range("D11").Value = 0.259491323550082
$ws.Range("E11").Value = 0.5094029088551439
$ws.Range("F11").Value = 0.3343990248860763
$ws.Range("G11").Value = 6
